$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 (sequence 1-6) with refreshed timestamps/content
$ws.Cells.Item(2, 1).Value = '2025-11-25T05:14:11.272Z'
$ws.Cells.Item(2, 10).Value = '{"Sequence":1,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"3aa492ca-003a-4aa8-a1eb-24189d6ec752","EventDtm":"2025-11-25T06:49:58Z","AppDtm":"2025-11-25T05:08:30Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T06:49:58Z","DeviceID":"JSGA622180057","GPSLockState":"UNLOCKED","GPSSatelliteCount":255,"GPSLastLock":3699,"GPSLatitude":17.244358,"GPSLongitude":78.443679,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"32291","CID":"244550678","RSSI":"-73","ExtPower":true,"ExtPowerVoltage":34.3,"BatteryVoltage":8.1,"DeviceTemp":30,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":null,"LastAssetRunState":null},"ReeferData":null}}'

$ws.Cells.Item(3, 1).Value = '2025-11-25T05:14:11.322Z'
$ws.Cells.Item(3, 10).Value = '{"Sequence":2,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"726959ca-0b44-4f97-9723-125c5a379a9e","EventDtm":"2025-11-25T05:09:08Z","AppDtm":"2025-11-25T05:09:30Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:09:08Z","DeviceID":"JSGA623040290","GPSLockState":"LOCKED","GPSSatelliteCount":16,"GPSLastLock":0,"GPSLatitude":17.661547,"GPSLongitude":83.089111,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"31121","CID":"232671233","RSSI":"-65","ExtPower":true,"ExtPowerVoltage":30,"BatteryVoltage":8,"DeviceTemp":34,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TDRU7153256","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T05:09:08Z","AssetType":"Reefer","AssetID":"TDRU7153256","OEM":"DAIKIN","TAmb":30.5,"TAmbQ":null,"TUSDA4":-72.8,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":50,"PSucQ":"asProvided","TDis":46.2,"TDisQ":null,"FreqComp":null,"TSuc":-14,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":1060,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":18,"PctExpVlvQ":null,"TEvap":3.5,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":75,"PctHumSetQ":"configured","FreqLine":50,"FreqLineQ":null,"VLine1":367.5,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":5.38,"TRtn1Q":null,"TRtn2":5.3,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"24C1","VerSWMinor":null,"PctSucVlv":98.78,"PctSucVlvQ":null,"TSup1":4.19,"TSup1Q":null,"TSup2":4.1,"TSup2Q":null,"AmpTotalDraw":14,"AmpTotalDrawQ":"asProvided","TUSDA1":-72.7,"TUSDA1Q":"OOR","TUSDA2":-72.8,"TUSDA2Q":"OOR","TUSDA3":-72.8,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"error","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-04-24T05:53:24Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'

$ws.Cells.Item(4, 1).Value = '2025-11-25T05:14:11.329Z'
$ws.Cells.Item(4, 10).Value = '{"Sequence":3,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"aecd3425-d975-4115-960d-d7fb04ad05bf","EventDtm":"2025-11-25T05:10:22Z","AppDtm":"2025-11-25T05:10:37Z","Events":["BatteryPowerOn"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:10:22Z","DeviceID":"JSGA623040193","GPSLockState":"LOCKED","GPSSatelliteCount":13,"GPSLastLock":0,"GPSLatitude":28.678715,"GPSLongitude":77.59958,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"97","LAC":"1827","CID":"230077975","RSSI":"-67","ExtPower":false,"ExtPowerVoltage":4.7,"BatteryVoltage":7.9,"DeviceTemp":26,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"SJKU4000017","LastAssetRunState":"Offline"},"ReeferData":null}}'

$ws.Cells.Item(5, 1).Value = '2025-11-25T05:14:11.338Z'
$ws.Cells.Item(5, 10).Value = '{"Sequence":4,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"30770747-3acc-4083-ae91-0ed242ad587d","EventDtm":"2025-11-25T05:11:08Z","AppDtm":"2025-11-25T05:11:20Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:11:08Z","DeviceID":"JSGA623040302","GPSLockState":"LOCKED","GPSSatelliteCount":10,"GPSLastLock":4,"GPSLatitude":19.252921,"GPSLongitude":73.016577,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"90","LAC":"6253","CID":"249578853","RSSI":"-65","ExtPower":true,"ExtPowerVoltage":28.9,"BatteryVoltage":8,"DeviceTemp":35,"RTDLOn":false,"VerFW":"W0206 1.41","DeviceType":"CT3500","DoorState":"Open","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TRIU6681882","LastAssetRunState":"Running"},"ReeferData":null}}'

$ws.Cells.Item(6, 1).Value = '2025-11-25T05:14:11.352Z'
$ws.Cells.Item(6, 10).Value = '{"Sequence":5,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"2b5837cc-f513-403f-b53a-83c0ced9ccc6","EventDtm":"2025-11-25T05:11:07Z","AppDtm":"2025-11-25T05:11:27Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:11:07Z","DeviceID":"JSGA623040329","GPSLockState":"UNLOCKED","GPSSatelliteCount":10,"GPSLastLock":12,"GPSLatitude":17.547937,"GPSLongitude":78.380464,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19392","CID":"232207115","RSSI":"-59","ExtPower":true,"ExtPowerVoltage":32.7,"BatteryVoltage":8,"DeviceTemp":28,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CCLU1035976","LastAssetRunState":"Offline"},"ReeferData":null}}'

$ws.Cells.Item(7, 1).Value = '2025-11-25T05:14:11.365Z'
$ws.Cells.Item(7, 10).Value = '{"Sequence":6,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"ed025364-a08b-41c1-936f-e5a3efb7e401","EventDtm":"2025-11-25T05:12:05Z","AppDtm":"2025-11-25T05:12:42Z","Events":["ACPowerOn"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:12:05Z","DeviceID":"JSGA622180045","GPSLockState":"LOCKED","GPSSatelliteCount":11,"GPSLastLock":0,"GPSLatitude":26.310615,"GPSLongitude":91.717618,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"56","LAC":"7134","CID":"250551313","RSSI":null,"ExtPower":true,"ExtPowerVoltage":30.6,"BatteryVoltage":8.1,"DeviceTemp":26,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TDRU7151905","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T05:12:05Z","AssetType":"Reefer","AssetID":"TDRU7151905","OEM":"DAIKIN","TAmb":23.19,"TAmbQ":null,"TUSDA4":-46.3,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":300,"PSucQ":"asProvided","TDis":22.9,"TDisQ":null,"FreqComp":null,"TSuc":20.9,"TSucQ":null,"MCond":"Off","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":300,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":5,"PctExpVlvQ":null,"TEvap":-0.06,"TEvapQ":null,"MCtrl3":"ThermoOff","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":75,"PctHumSetQ":"configured","FreqLine":50,"FreqLineQ":null,"VLine1":404.4,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"FanReduction","TRtn1":46.62,"TRtn1Q":null,"TRtn2":11.2,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"24C2","VerSWMinor":null,"PctSucVlv":50,"PctSucVlvQ":null,"TSup1":11,"TSup1Q":null,"TSup2":10.8,"TSup2Q":null,"AmpTotalDraw":0,"AmpTotalDrawQ":"asProvided","TUSDA1":-46.3,"TUSDA1Q":"OOR","TUSDA2":-46.3,"TUSDA2Q":"OOR","TUSDA3":-46.3,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"OOR","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":null,"PTIResult":null,"TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'

$ws.Cells.Item(8, 1).Value = '2025-11-25T05:14:11.370Z'
$ws.Cells.Item(8, 10).Value = '{"Sequence":7,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"55df62f7-09ae-4619-8cff-e36059733f52","EventDtm":"2025-11-25T05:13:03Z","AppDtm":"2025-11-25T05:13:22Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:13:03Z","DeviceID":"JSGA623040284","GPSLockState":"LOCKED","GPSSatelliteCount":15,"GPSLastLock":0,"GPSLatitude":17.537149,"GPSLongitude":78.175627,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19397","CID":"247361559","RSSI":"-73","ExtPower":true,"ExtPowerVoltage":27.6,"BatteryVoltage":8,"DeviceTemp":35,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TRIU6618637","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T05:13:03Z","AssetType":"Reefer","AssetID":"TRIU6618637","OEM":"CARRIER","TAmb":28.09,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":0,"PctCO2Q":null,"PctCO2Set":5,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":0,"PCond":174.4,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04898885","AmpPhA":8.93,"AmpPhB":10.39,"AmpPhC":8.93,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":-23.8,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":0,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":11.25,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":394.61,"VLine2":null,"VLine3":null,"MEvapFanLS":1,"PctO2":0,"PctO2Q":"asProvided","PctO2Set":10,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":-17.83,"TRtn1Q":null,"TRtn2":-17.87,"TRtn2Q":null,"TSet":-20,"TSetQ":null,"VerSWMajor":"5168","VerSWMinor":null,"PctSucVlv":100,"PctSucVlvQ":null,"TSup1":-32.2,"TSup1Q":null,"TSup2":-31.59,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-09-30T10:27:06Z","PTIResult":"Skipped","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":[{"OemAlarm":52,"RCAlias":"AL52","RCSeverity":"Minor","Active":true},{"OemAlarm":53,"RCAlias":"AL53","RCSeverity":"Informational","Active":true}]}}}'
$ws.Cells.Item(8, 2).Value = "N/A" ; $ws.Cells.Item(8, 3).Value = "N/A" ; $ws.Cells.Item(8, 4).Value = "N/A" ; $ws.Cells.Item(8, 5).Value = "N/A" ; $ws.Cells.Item(8, 6).Value = "N/A" ; $ws.Cells.Item(8, 7).Value = "N/A" ; $ws.Cells.Item(8, 8).Value = "N/A" ; $ws.Cells.Item(8, 9).Value = "N/A"

$ws.Cells.Item(9, 1).Value = '2025-11-25T05:14:11.475Z'
$ws.Cells.Item(9, 10).Value = '{"Sequence":8,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"e0ae646c-e282-4094-80ed-29300cf60aba","EventDtm":"2025-11-25T05:13:32Z","AppDtm":"2025-11-25T05:13:49Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:13:32Z","DeviceID":"JSGA622340343","GPSLockState":"LOCKED","GPSSatelliteCount":14,"GPSLastLock":0,"GPSLatitude":12.804881,"GPSLongitude":77.662177,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9003","CID":"46904076","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.4,"BatteryVoltage":8,"DeviceTemp":32,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CCLU1035976","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T05:13:32Z","AssetType":"Reefer","AssetID":"CCLU1035976","OEM":"CARRIER","TAmb":32.44,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":5,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":101.88,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04904656","AmpPhA":11.25,"AmpPhB":10.76,"AmpPhC":11.25,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":6.21,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":412.62,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":-0.03,"PctO2Q":"OOR","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.8,"TRtn1Q":null,"TRtn2":5.81,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"5180","VerSWMinor":null,"PctSucVlv":3.71,"PctSucVlvQ":null,"TSup1":4.97,"TSup1Q":null,"TSup2":4.92,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-05-12T10:01:32Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'
$ws.Cells.Item(9, 2).Value = "N/A" ; $ws.Cells.Item(9, 3).Value = "N/A" ; $ws.Cells.Item(9, 4).Value = "N/A" ; $ws.Cells.Item(9, 5).Value = "N/A" ; $ws.Cells.Item(9, 6).Value = "N/A" ; $ws.Cells.Item(9, 7).Value = "N/A" ; $ws.Cells.Item(9, 8).Value = "N/A" ; $ws.Cells.Item(9, 9).Value = "N/A"

$ws.Cells.Item(10, 1).Value = '2025-11-25T05:15:00.123Z'
$ws.Cells.Item(10, 10).Value = '{"Sequence":9,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"6d975b75-d8a9-4d89-81fe-d051c2f61815","EventDtm":"2025-11-25T05:14:43Z","AppDtm":"2025-11-25T05:14:59Z","Events":["MicroAlarm"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:14:43Z","DeviceID":"JSGA622180045","GPSLockState":"LOCKED","GPSSatelliteCount":12,"GPSLastLock":0,"GPSLatitude":26.310582,"GPSLongitude":91.717582,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"56","LAC":"7134","CID":"250551313","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.5,"BatteryVoltage":8.1,"DeviceTemp":28,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TDRU7151905","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T05:14:43Z","AssetType":"Reefer","AssetID":"TDRU7151905","OEM":"DAIKIN","TAmb":23.38,"TAmbQ":null,"TUSDA4":-46.3,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":310,"PSucQ":"asProvided","TDis":22.9,"TDisQ":null,"FreqComp":null,"TSuc":21.1,"TSucQ":null,"MCond":"Off","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":300,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":50,"PctExpVlvQ":null,"TEvap":-0.06,"TEvapQ":null,"MCtrl3":"PTIEmergencyStop1","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":75,"PctHumSetQ":"configured","FreqLine":50,"FreqLineQ":null,"VLine1":397.8,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"FanReduction","TRtn1":46.38,"TRtn1Q":null,"TRtn2":10.6,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"24C2","VerSWMinor":null,"PctSucVlv":50.3,"PctSucVlvQ":null,"TSup1":9.88,"TSup1Q":null,"TSup2":9.8,"TSup2Q":null,"AmpTotalDraw":2,"AmpTotalDrawQ":"asProvided","TUSDA1":-46.3,"TUSDA1Q":"OOR","TUSDA2":-46.3,"TUSDA2Q":"OOR","TUSDA3":-46.3,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"OOR","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":null,"PTIResult":null,"TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":[{"OemAlarm":409,"RCAlias":"E409","RCSeverity":"Informational","Active":true}]}}}'
$ws.Cells.Item(10, 2).Value = "N/A" ; $ws.Cells.Item(10, 3).Value = "N/A" ; $ws.Cells.Item(10, 4).Value = "N/A" ; $ws.Cells.Item(10, 5).Value = "N/A" ; $ws.Cells.Item(10, 6).Value = "N/A" ; $ws.Cells.Item(10, 7).Value = "N/A" ; $ws.Cells.Item(10, 8).Value = "N/A" ; $ws.Cells.Item(10, 9).Value = "N/A"

$ws.Cells.Item(11, 1).Value = '2025-11-25T05:16:10.940Z'
$ws.Cells.Item(11, 10).Value = '{"Sequence":10,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"0f08d980-a463-47ed-821b-34277681d5a9","EventDtm":"2025-11-25T05:15:54Z","AppDtm":"2025-11-25T05:16:10Z","Events":["MicroAlarm"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T05:15:54Z","DeviceID":"JSGA622180045","GPSLockState":"LOCKED","GPSSatelliteCount":14,"GPSLastLock":0,"GPSLatitude":26.3106,"GPSLongitude":91.717571,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"56","LAC":"7134","CID":"250551313","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.1,"BatteryVoltage":8.1,"DeviceTemp":28,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TDRU7151905","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T05:15:54Z","AssetType":"Reefer","AssetID":"TDRU7151905","OEM":"DAIKIN","TAmb":23.38,"TAmbQ":null,"TUSDA4":-46.4,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":180,"PSucQ":"asProvided","TDis":22.9,"TDisQ":null,"FreqComp":null,"TSuc":18.2,"TSucQ":null,"MCond":"Off","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":1090,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":100,"PctExpVlvQ":null,"TEvap":-0.06,"TEvapQ":null,"MCtrl3":"FullCool","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"Off","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":75,"PctHumSetQ":"configured","FreqLine":50,"FreqLineQ":null,"VLine1":399.8,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"CoolMax","TRtn1":-0.06,"TRtn1Q":null,"TRtn2":11.2,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"24C2","VerSWMinor":null,"PctSucVlv":95.43,"PctSucVlvQ":null,"TSup1":10.69,"TSup1Q":null,"TSup2":10.5,"TSup2Q":null,"AmpTotalDraw":0,"AmpTotalDrawQ":"asProvided","TUSDA1":-46.3,"TUSDA1Q":"OOR","TUSDA2":-46.3,"TUSDA2Q":"OOR","TUSDA3":-46.3,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"OOR","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":null,"PTIResult":null,"TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":[{"OemAlarm":403,"RCAlias":"E403","RCSeverity":"Minor","Active":true},{"OemAlarm":409,"RCAlias":"E409","RCSeverity":"Informational","Active":true}]}}}'
$ws.Cells.Item(11, 2).Value = "N/A" ; $ws.Cells.Item(11, 3).Value = "N/A" ; $ws.Cells.Item(11, 4).Value = "N/A" ; $ws.Cells.Item(11, 5).Value = "N/A" ; $ws.Cells.Item(11, 6).Value = "N/A" ; $ws.Cells.Item(11, 7).Value = "N/A" ; $ws.Cells.Item(11, 8).Value = "N/A" ; $ws.Cells.Item(11, 9).Value = "N/A"
